# Document.ContainsFields.docx edit script
#
# Ports the document's <w:fldSimple> fields to "expanded" field-code form
# (begin/instrText/separate/result/end run sequences), refreshes the
# DATE/TIME field results, adds a noProof paragraph mark to the Author
# paragraph, and appends a new "Field with no separator:" paragraph
# containing a GOTOBUTTON field (begin/instrText/end - no "separate"
# fldChar) that now also carries the trailing _GoBack bookmark.

$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1) DATE field result text: refresh in place (keeps the existing
#    begin/instrText/separate/end run structure - only the visible result
#    changes).
$d.Content.Find.Execute("Tuesday, 21 August 2018", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Monday, 28 January 2019", 2) | Out-Null

# 2) "Time:" paragraph - expand the <w:fldSimple> TIME field into explicit
#    fldChar begin/instrText/separate/text/end runs, with the new result.
$timePara = $d.Paragraphs.Item(2)
$timeXml = "<w:p $w>" + `
    '<w:r><w:t xml:space="preserve">Time: </w:t></w:r>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:instrText xml:space="preserve"> TIME   \* MERGEFORMAT </w:instrText></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:t>12:06 AM</w:t></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>' + `
    '</w:p>'
$timePara.Range.InsertXML($timeXml) | Out-Null

# 3) "Filename:" paragraph - expand the FILENAME fldSimple field.
$filenamePara = $d.Paragraphs.Item(3)
$filenameXml = "<w:p $w>" + `
    '<w:r><w:t>Filename:</w:t></w:r>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:instrText xml:space="preserve"> FILENAME   \* MERGEFORMAT </w:instrText></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:t>Document.ContainsFields.docx</w:t></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>' + `
    '</w:p>'
$filenamePara.Range.InsertXML($filenameXml) | Out-Null

# 4) "Size:" paragraph - expand the FILESIZE fldSimple field, keep the
#    trailing " bytes" runs.
$sizePara = $d.Paragraphs.Item(4)
$sizeXml = "<w:p $w>" + `
    '<w:r><w:t>Size:</w:t></w:r>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:instrText xml:space="preserve"> FILESIZE   \* MERGEFORMAT </w:instrText></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:t>11632</w:t></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>bytes</w:t></w:r>' + `
    '</w:p>'
$sizePara.Range.InsertXML($sizeXml) | Out-Null

# 5) "Author:" paragraph - give the paragraph mark a noProof rPr, expand
#    the AUTHOR fldSimple field, and (in the same InsertXML call) append a
#    brand-new paragraph right after it that demonstrates a field with no
#    "separate" fldChar (GOTOBUTTON). The trailing _GoBack bookmark that
#    used to sit at the end of the Author paragraph now moves to the end
#    of this new paragraph.
$authorPara = $d.Paragraphs.Item(5)
$authorAndGotoXml = "<w:p $w>" + `
    '<w:pPr><w:rPr><w:noProof/></w:rPr></w:pPr>' + `
    '<w:r><w:t>Author:</w:t></w:r>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:instrText xml:space="preserve"> AUTHOR  \* FirstCap  \* MERGEFORMAT </w:instrText></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:t>Roman S</w:t></w:r>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>' + `
    '</w:p>' + `
    "<w:p $w>" + `
    '<w:pPr><w:tabs><w:tab w:val="center" w:pos="4513"/></w:tabs></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Field with no separator: </w:t></w:r>' + `
    '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' + `
    '<w:r><w:instrText xml:space="preserve"> GOTOBUTTON  </w:instrText></w:r>' + `
    '<w:r><w:fldChar w:fldCharType="end"/></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '</w:p>'
$authorPara.Range.InsertXML($authorAndGotoXml) | Out-Null
